$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'296.97"
$ws.Range("E2").Value = "'1.68%"

$ws.Range("D3").Value = "'41.82"
$ws.Range("E3").Value = "'3.82%"

$ws.Range("D4").Value = "'5.004"
$ws.Range("E4").Value = "'-0.02%"

$ws.Range("D5").Value = "'0.07516"
$ws.Range("E5").Value = "'2.73%"

$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.381"
$ws.Range("E6").Value = "'2.01%"

$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.583"
$ws.Range("E7").Value = "'4.13%"

$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9264"
$ws.Range("E8").Value = "'-0.12%"

$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.401"
$ws.Range("E9").Value = "'0.97%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1199"
$ws.Range("E10").Value = "'0.08%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1830"
$ws.Range("E11").Value = "'5.09%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08869"
$ws.Range("E12").Value = "'2.73%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04073"
$ws.Range("E13").Value = "'-5.95%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1048"
$ws.Range("E14").Value = "'-0.61%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001284"
$ws.Range("E15").Value = "'1.00%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005943"
$ws.Range("E16").Value = "'-0.25%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.357"
$ws.Range("E17").Value = "'0.55%"

$ws.Range("D19").Value = "'8.129"
$ws.Range("E19").Value = "'1.93%"

$ws.Range("D20").Value = "'0.1390"
$ws.Range("E20").Value = "'-0.02%"

$ws.Range("E21").Value = "'11.00%"

$ws.Range("D22").Value = "'0.04101"
$ws.Range("E22").Value = "'4.30%"

$ws.Range("D23").Value = "'0.001265"
$ws.Range("E23").Value = "'0.32%"

$ws.Range("D24").Value = "'0.003906"
$ws.Range("E24").Value = "'3.40%"

$ws.Range("E25").Value = "'-3.96%"

$ws.Range("D38").Value = "'0.02405"
$ws.Range("E38").Value = "'5.20%"

$ws.Range("D39").Value = "'0.05239"
$ws.Range("E39").Value = "'5.20%"

$ws.Range("D40").Value = "'0.006303"
$ws.Range("E40").Value = "'17.66%"

$ws.Range("D41").Value = "'0.007816"
$ws.Range("E41").Value = "'1.49%"

$ws.Range("D42").Value = "'0.1326"
$ws.Range("E42").Value = "'3.37%"

$ws.Range("D43").Value = "'0.007390"
$ws.Range("E43").Value = "'0.85%"

$ws.Range("D44").Value = "'0.007810"
$ws.Range("E44").Value = "'-1.03%"

$ws.Range("D45").Value = "'0.2960"
$ws.Range("E45").Value = "'-6.64%"

$ws.Range("D46").Value = "'0.00006527"
$ws.Range("E46").Value = "'3.24%"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.08%"

$ws.Range("D48").Value = "'0.03143"
$ws.Range("E48").Value = "'53.90%"

$ws.Range("D49").Value = "'0.004202"
$ws.Range("E49").Value = "'0.01%"

$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.08%"

$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.08%"
